$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Ensure the touched cells keep a text number format so values are written
# back as strings (matching the original t="str" cell type) rather than
# being reinterpreted as numeric.
$cells = "C2","D2","E2","F2","C4","D4","C5","D5","E5","C6","D6","F6"
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("C2").Value = "25"
$ws.Range("D2").Value = "8"
$ws.Range("E2").Value = "1"
$ws.Range("F2").Value = "3"

# Row 4
$ws.Range("C4").Value = "0"
$ws.Range("D4").Value = "2"

# Row 5
$ws.Range("C5").Value = "3"
$ws.Range("D5").Value = "4"
$ws.Range("E5").Value = "0"

# Row 6
$ws.Range("C6").Value = "4"
$ws.Range("D6").Value = "2"
$ws.Range("F6").Value = "0"
